$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new "through" date.
$ws.Name = "Through 2022-02-20"

# Update the February row label text.
$ws.Range("A3").Value = "February (through 02-20)"

# Update February row (row 3) values.
$ws.Range("B3").Value = 8
$ws.Range("C3").Value = 26
$ws.Range("F3").Value = 21
$ws.Range("G3").Value = 50
$ws.Range("H3").Value = 88
$ws.Range("I3").Value = 98

# Update Total row (row 4) values.
$ws.Range("B4").Value = 34
$ws.Range("C4").Value = 77
$ws.Range("F4").Value = 70
$ws.Range("G4").Value = 124
$ws.Range("H4").Value = 305
$ws.Range("I4").Value = 257
